# 自动更新价格数据：插入新的最新日期行，原有数据下移一行
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 在第2行前插入一行，原第2行（2025-11-21 的数据）整体下移到第3行
$ws.Rows("2:2").Insert()

# 在新的第2行写入最新一天的数据
# 日期以文本形式写入（前导单引号强制为文本，避免被识别为日期序列值）
$ws.Range("A2").Value = "'2025-11-22"
$ws.Range("B2").Value = 783.5
$ws.Range("C2").Value = 1112
$ws.Range("D2").Value = 3610

# 清除插入行从上方继承的格式，使其与原始数据行保持一致（无额外样式）
$ws.Rows("2:2").ClearFormats()
